# Weekly update: insert 3 new rows (new week 2022-06-02 / serial 44714)
# at the top of the existing "Vega Monumental Concepción - Kiwi" block
# (rows 140-142), pushing the previous rows 140-144 down to 143-147.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before the old row 140; this shifts old rows
# 140-144 down to 143-147, preserving their content and formatting.
$ws.Rows("140:142").Insert()

# Shared values for this Mercado/Producto block (identical across all
# rows in this subset of the sheet).
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$fecha     = 44714
$codreg    = 8
$tipo      = "Fruta"
$prodId    = 100101
$producto  = "Berries"
$catId     = 100101007
$categoria = "Kiwi"
$variedad  = "Hayward"
$unidad    = '$/bandeja 18 kilos'
$origen    = "Región de O'Higgins"
$kgUnidad  = 18

# New row 140: Especial, volumen 50, 11000/11000/11000, 611 $/Kg
$r = 140
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 50
$ws.Cells.Item($r, 14).Value = 11000
$ws.Cells.Item($r, 15).Value = 11000
$ws.Cells.Item($r, 16).Value = 11000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 611
$ws.Cells.Item($r, 20).Value = $kgUnidad

# New row 141: Primera, volumen 50, 9000/9000/9000, 500 $/Kg
$r = 141
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 50
$ws.Cells.Item($r, 14).Value = 9000
$ws.Cells.Item($r, 15).Value = 9000
$ws.Cells.Item($r, 16).Value = 9000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 500
$ws.Cells.Item($r, 20).Value = $kgUnidad

# New row 142: Segunda, volumen 50, 8000/8000/8000, 444 $/Kg
$r = 142
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $prodId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $catId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 50
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 8000
$ws.Cells.Item($r, 16).Value = 8000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 444
$ws.Cells.Item($r, 20).Value = $kgUnidad
